# Updates for 4/16 status meeting
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 (R20 risk) status changes from Open to Closed, and gets the
# "closed" grey fill style applied (same style used for other closed rows).
$ws.Range("E21").Value = "Closed"
$ws.Range("A21:F21").Interior.ThemeColor = 1
$ws.Range("A21:F21").Interior.TintAndShade = -0.249977111117893

# New row 22: R21 risk about transition from Rashmi to Mike Hunter.
$ws.Range("A22").Value = "R21"
$ws.Range("B22").Value = "Transition from Rashmi to Mike Hunter will result in loss of historical knowledge and project context"
$ws.Range("C22").Value = "Ulli, Juli and JJ"
$ws.Range("D22").Value = "High"
$ws.Range("E22").Value = "Open"
$ws.Range("F22").Value = "5AM submitted a training plan for Mike's first two months to mitigate the impact of Rashmi leaving the project.  The plan was accepted and is underway.  Mike also has access to Will Fitzhugh and Todd Parnell as a recourse if significant issues arise."

$ws.Range("A22:F22").Font.Name = "Times New Roman"
$ws.Range("A22:F22").Font.Size = 12
$ws.Range("A22:F22").VerticalAlignment = -4108
$ws.Range("A22:F22").WrapText = $true

$ws.Rows.Item(21).RowHeight = 90
$ws.Rows.Item(22).RowHeight = 75

$ws.Range("A1:F22").Select
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("F23").Select
